$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.603.10"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.89%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.880.73"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.03%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.29%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.76"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.22%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.77%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5111"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.44%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3942"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.81%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08426"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.14%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.116"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.52%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.73"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.46%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.292"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.09%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.880.61"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.49%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.53"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.46%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.291"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.88%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.009"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.24%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001109"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.79%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.58"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.81%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06729"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.02%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.77"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.26%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.007"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.54%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.985"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.27%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.631.31"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.89%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.16"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.48%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.254"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.38%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.096.22"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.51%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.16"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.23%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.81"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.34%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.383"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.43%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.87"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.74%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1055"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.81%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.058"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.03%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.842"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.52%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.617"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.74%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02463"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.24%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06552"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.09%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2191"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.18%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.974"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.73%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.266"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.23%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.68%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6497"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.25%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.096"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.71%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.22"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.01%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.21%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6087"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.02%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.04"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.39%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.702"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.49%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.042"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.06%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.221"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.67%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "122.75"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.34%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -7.75%  "
